# "Generate Report for Handoff"
# Updates the localization-status report: the b.md row moves from
# "Handed back: in sync with en-US" to "Ready for handoff" on the Overview
# sheet, and on each locale sheet (zh-cn, de-de) a fresh handoff file /
# datetime / error-detail get recorded for b.md. Also widens the
# "Error Detail" column (P) on both locale sheets so the longer message fits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet — row for b.md (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-27 16:37:15"

# ---------------------------------------------------------------------
# zh-cn sheet — row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-27 16:37:10"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b87d76c01eada3a59ce85294e001b79df8629ab/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ade5e2d1dfaf516056d46e5910af22f1a00370a/e2e/b.md."

# widen the Error Detail column (P) so the long message is readable
$zhcn.Range("P1").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet — row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-27 16:37:15"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b87d76c01eada3a59ce85294e001b79df8629ab/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ade5e2d1dfaf516056d46e5910af22f1a00370a/e2e/b.md."

# widen the Error Detail column (P) so the long message is readable
$dede.Range("P1").ColumnWidth = 39.166666666666664
